# Update TPM-derived NATMI metrics (Receptor / Edge expression & specificity
# columns) on Sheet1 to reflect the re-run of the scripts with new TPM values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("M2").Value = 3.867218333333334
$ws.Range("N2").Value = 11.601655
$ws.Range("O2").Value = 0.1566152977872902
$ws.Range("P2").Value = 0.1566152977872902
$ws.Range("Q2").Value = 0.2128452517027778
$ws.Range("R2").Value = 1.915607265325
$ws.Range("S2").Value = 0.1566152977872902
$ws.Range("T2").Value = 0.1566152977872902

# Row 3
$ws.Range("N3").Value = 33.813685
$ws.Range("O3").Value = 0.4564642152831324
$ws.Range("P3").Value = 0.4564642152831324
$ws.Range("Q3").Value = 0.6203496220861111
$ws.Range("S3").Value = 0.4564642152831324
$ws.Range("T3").Value = 0.4564642152831324

# Row 4
$ws.Range("M4").Value = 5.654344666666667
$ws.Range("N4").Value = 16.963034
$ws.Range("O4").Value = 0.2289906587711778
$ws.Range("P4").Value = 0.2289906587711778
$ws.Range("Q4").Value = 0.3112057065455556
$ws.Range("R4").Value = 2.80085135891
$ws.Range("S4").Value = 0.2289906587711778
$ws.Range("T4").Value = 0.2289906587711778

# Row 5
$ws.Range("M5").Value = 0.819389
$ws.Range("N5").Value = 2.458167
$ws.Range("O5").Value = 0.03318376186120772
$ws.Range("P5").Value = 0.03318376186120772
$ws.Range("Q5").Value = 0.04509780491166667
$ws.Range("R5").Value = 0.405880244205
$ws.Range("S5").Value = 0.03318376186120772
$ws.Range("T5").Value = 0.03318376186120772

# Row 6
$ws.Range("M6").Value = 3.080288333333333
$ws.Range("N6").Value = 9.240864999999999
$ws.Range("O6").Value = 0.1247460662971919
$ws.Range("P6").Value = 0.1247460662971919
$ws.Range("Q6").Value = 0.1695339360527778
$ws.Range("R6").Value = 1.525805424475
$ws.Range("S6").Value = 0.1247460662971919
$ws.Range("T6").Value = 0.1247460662971919
